$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 for the new "08/2025" month entry (day 1),
# pushing all existing data rows (previously 2..66) down by one (3..67).
$ws.Rows.Item(2).Insert()

# The inserted row copies formatting from the row below it (bold/border style);
# the source data rows carry no explicit style, so clear it to match.
$ws.Range("A2:E2").ClearFormats()

# New August data point (Dia=1, total_venda=23013.08, Mes=8, Ano=2025, Periodo=08/2025)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 23013.08
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 2025
$ws.Range("E2").Value = "08/2025"

# Refreshed total_venda values for a handful of July days (rows shifted down by
# one after the insert above: old row N is now row N+1).
$ws.Range("B3").Value = 18112.93
$ws.Range("B12").Value = 18544.3
$ws.Range("B16").Value = 7505.85
$ws.Range("B17").Value = 499270.27
$ws.Range("B22").Value = 18613.59
$ws.Range("B24").Value = 59573.83
$ws.Range("B25").Value = 25281
